$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Row 4: new work-hours entry.
# A4: date 2020-08-04, same numeric-date style ("d-mmm", s=1) as A2/A3.
$ws.Range("A4").Value2 = 44047
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat

# B4: start time 15:30, same time style ("h:mm", s=2) as B2/B3.
$ws.Range("B4").Value2 = 0.64583333333333337
$ws.Range("B4").NumberFormat = $ws.Range("B3").NumberFormat

# C4: end time not known yet -> literal "?" text (keeps default/general style).
$ws.Range("C4").Value = "?"

# D4: duration not known yet -> "?" (keeps existing right-aligned style, s=3).
$ws.Range("D4").Value = "?"

# E4: description, reuse text from row 3 ("tietokannan maarittelya").
$ws.Range("E4").Value = $ws.Range("E3").Value2

# Move the active selection to A5, matching the saved cursor position.
$ws.Range("A5").Select() | Out-Null
